$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 9588.235000000001
$ws.Range("I76").Value = 11500
$ws.Range("K76").Value = 11500
$ws.Range("M76").Value = -11185
$ws.Range("H79").Value = 9588.235000000001
$ws.Range("I79").Value = 11500
$ws.Range("K79").Value = 11500
$ws.Range("M79").Value = -10408
$ws.Range("H137").Value = 4413251
$ws.Range("I137").Value = 2084568
$ws.Range("J137").Value = 10002090
$ws.Range("K137").Value = 6253704
$ws.Range("L137").Value = 30006270
$ws.Range("M137").Value = -6251154
$ws.Range("N137").Value = -30011370
$ws.Range("H141").Value = 3070.6667
$ws.Range("I141").Value = 2481.8333
$ws.Range("J141").Value = 3659.5
$ws.Range("K141").Value = 7445.499899999999
$ws.Range("L141").Value = 10978.5
$ws.Range("M141").Value = -2265.499899999999
$ws.Range("N141").Value = -21338.5

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3435.7827
$ws.Range("I2").Value = 1109.5333
$ws.Range("J2").Value = 7797.5
$ws.Range("K2").Value = 1109.5333
$ws.Range("L2").Value = 7797.5
$ws.Range("M2").Value = -996.5333000000001
$ws.Range("N2").Value = -8023.5
$ws.Range("H59").Value = 29877
$ws.Range("J59").Value = 29877
$ws.Range("L59").Value = 29877
$ws.Range("N59").Value = -31485
$ws.Range("H61").Value = 1174.5
$ws.Range("I61").Value = 1054.6842
$ws.Range("J61").Value = 1933.3334
$ws.Range("K61").Value = 1054.6842
$ws.Range("L61").Value = 1933.3334
$ws.Range("M61").Value = -842.6841999999999
$ws.Range("N61").Value = -2357.3334
$ws.Range("H74").Value = 695.67444
$ws.Range("I74").Value = 671.6842
$ws.Range("J74").Value = 878
$ws.Range("K74").Value = 671.6842
$ws.Range("L74").Value = 878
$ws.Range("M74").Value = 202.3158
$ws.Range("N74").Value = -2626
$ws.Range("H77").Value = 695.67444
$ws.Range("I77").Value = 671.6842
$ws.Range("J77").Value = 878
$ws.Range("K77").Value = 3358.421
$ws.Range("L77").Value = 4390
$ws.Range("M77").Value = 1009.579
$ws.Range("N77").Value = -13126
$ws.Range("H110").Value = 3234.0256
$ws.Range("I110").Value = 2586.5
$ws.Range("J110").Value = 11004.333
$ws.Range("K110").Value = 2586.5
$ws.Range("L110").Value = 11004.333
$ws.Range("M110").Value = -541.5
$ws.Range("N110").Value = -15094.333
$ws.Range("H116").Value = 3435.7827
$ws.Range("I116").Value = 1109.5333
$ws.Range("J116").Value = 7797.5
$ws.Range("K116").Value = 1109.5333
$ws.Range("L116").Value = 7797.5
$ws.Range("M116").Value = 1184.4667
$ws.Range("N116").Value = -12385.5
$ws.Range("H132").Value = 298551.66
$ws.Range("I132").Value = 456851.47
$ws.Range("J132").Value = 8335.333000000001
$ws.Range("K132").Value = 1370554.41
$ws.Range("L132").Value = 25005.999
$ws.Range("M132").Value = -1368024.41
$ws.Range("N132").Value = -30065.999
$ws.Range("H136").Value = 1174.5
$ws.Range("I136").Value = 1054.6842
$ws.Range("J136").Value = 1933.3334
$ws.Range("K136").Value = 3164.0526
$ws.Range("L136").Value = 5800.0002
$ws.Range("M136").Value = -614.0526
$ws.Range("N136").Value = -10900.0002

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3435.7827
$ws.Range("I3").Value = 1109.5333
$ws.Range("J3").Value = 7797.5
$ws.Range("K3").Value = 1109.5333
$ws.Range("L3").Value = 7797.5
$ws.Range("M3").Value = -995.5333000000001
$ws.Range("N3").Value = -8025.5
$ws.Range("H134").Value = 117503.46
$ws.Range("I134").Value = 127091.25
$ws.Range("J134").Value = 2450
$ws.Range("K134").Value = 381273.75
$ws.Range("L134").Value = 7350
$ws.Range("M134").Value = -378738.75
$ws.Range("N134").Value = -12420

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 300
$ws.Range("I2").Value = 300
$ws.Range("K2").Value = 300
$ws.Range("M2").Value = -187
$ws.Range("H22").Value = 537.5
$ws.Range("I22").Value = 550
$ws.Range("J22").Value = 475
$ws.Range("K22").Value = 550
$ws.Range("L22").Value = 475
$ws.Range("M22").Value = -200
$ws.Range("N22").Value = -1175
$ws.Range("H31").Value = 1903.1346
$ws.Range("I31").Value = 1621.1951
$ws.Range("J31").Value = 2954
$ws.Range("K31").Value = 1621.1951
$ws.Range("L31").Value = 2954
$ws.Range("M31").Value = -1326.1951
$ws.Range("N31").Value = -3544
$ws.Range("H34").Value = 1903.1346
$ws.Range("I34").Value = 1621.1951
$ws.Range("J34").Value = 2954
$ws.Range("K34").Value = 1621.1951
$ws.Range("L34").Value = 2954
$ws.Range("M34").Value = -1419.1951
$ws.Range("N34").Value = -3358
$ws.Range("H35").Value = 8531.532999999999
$ws.Range("I35").Value = 1734.375
$ws.Range("J35").Value = 16299.714
$ws.Range("K35").Value = 1734.375
$ws.Range("L35").Value = 16299.714
$ws.Range("M35").Value = -1440.375
$ws.Range("N35").Value = -16887.714
$ws.Range("H58").Value = 1472.4286
$ws.Range("I58").Value = 1524.5
$ws.Range("J58").Value = 1160
$ws.Range("K58").Value = 1524.5
$ws.Range("L58").Value = 1160
$ws.Range("M58").Value = -1321.5
$ws.Range("N58").Value = -1566
$ws.Range("H132").Value = 2897.739
$ws.Range("I132").Value = 2255.7334
$ws.Range("J132").Value = 4101.5
$ws.Range("K132").Value = 6767.2002
$ws.Range("L132").Value = 12304.5
$ws.Range("M132").Value = -4237.2002
$ws.Range("N132").Value = -17364.5
$ws.Range("H134").Value = 5229.069
$ws.Range("I134").Value = 5745.16
$ws.Range("J134").Value = 2003.5
$ws.Range("K134").Value = 17235.48
$ws.Range("L134").Value = 6010.5
$ws.Range("M134").Value = -14700.48
$ws.Range("N134").Value = -11080.5
$ws.Range("H136").Value = 1472.4286
$ws.Range("I136").Value = 1524.5
$ws.Range("J136").Value = 1160
$ws.Range("K136").Value = 4573.5
$ws.Range("L136").Value = 3480
$ws.Range("M136").Value = -2023.5
$ws.Range("N136").Value = -8580

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 172.5
$ws.Range("I4").Value = 96
$ws.Range("K4").Value = 288
$ws.Range("M4").Value = -176
$ws.Range("H49").Value = 6548.8
$ws.Range("J49").Value = 6943.1113
$ws.Range("L49").Value = 20829.3339
$ws.Range("N49").Value = -21141.3339

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1776.3
$ws.Range("I102").Value = 1621.8
$ws.Range("J102").Value = 2239.8
$ws.Range("K102").Value = 1621.8
$ws.Range("L102").Value = 2239.8
$ws.Range("M102").Value = 0.2000000000000455
$ws.Range("N102").Value = -5483.8
$ws.Range("H132").Value = 2646.44
$ws.Range("I132").Value = 2128.2
$ws.Range("J132").Value = 4719.4
$ws.Range("K132").Value = 6384.599999999999
$ws.Range("L132").Value = 14158.2
$ws.Range("M132").Value = -3854.599999999999
$ws.Range("N132").Value = -19218.2

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1000
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H35").Value = 767.9
$ws.Range("I35").Value = 767.9
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 767.9
$ws.Range("L35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -431.9
$ws.Range("H132").Value = 2225
$ws.Range("I132").Value = 1342.909
$ws.Range("K132").Value = 4028.727
$ws.Range("M132").Value = -1498.727
$ws.Range("H136").Value = 1801.0834
$ws.Range("I136").Value = 1582.875
$ws.Range("J136").Value = 2237.5
$ws.Range("K136").Value = 4748.625
$ws.Range("L136").Value = 6712.5
$ws.Range("M136").Value = -2198.625
$ws.Range("N136").Value = -11812.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2496
$ws.Range("I126").Value = 1883.5555
$ws.Range("J126").Value = 4333.3335
$ws.Range("K126").Value = 5650.666499999999
$ws.Range("L126").Value = 13000.0005
$ws.Range("M126").Value = -3180.666499999999
$ws.Range("N126").Value = -17940.0005
$ws.Range("H132").Value = 2608.6758
$ws.Range("I132").Value = 2210.3914
$ws.Range("J132").Value = 3263
$ws.Range("K132").Value = 6631.174199999999
$ws.Range("L132").Value = 9789
$ws.Range("M132").Value = -4101.174199999999
$ws.Range("N132").Value = -14849
$ws.Range("H136").Value = 1154.8379
$ws.Range("I136").Value = 964.1142599999999
$ws.Range("J136").Value = 4492.5
$ws.Range("K136").Value = 2892.34278
$ws.Range("L136").Value = 13477.5
$ws.Range("M136").Value = -342.3427799999999
$ws.Range("N136").Value = -18577.5
